$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1), matching the formatting of the other
# header cells (bold / bordered / centered style already used in G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats - reuse existing style, no new one

# New data values for the Save column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
